$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new "Good Polygon" values on rows that were previously blank
$ws.Range("G2").Value = "Good Polygon"
$ws.Range("G4").Value = "Good Polygon"
$ws.Range("G13").Value = "Good Polygon"

# Clear out rows that previously had a validation value
$ws.Range("G5").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("G16").Value = ""
$ws.Range("G17").Value = ""
$ws.Range("G19").Value = ""
